# Fruta / hortaliza, semanal
# Insert a new weekly record (3 rows) at the top of the price-history block
# (rows 880-882), pushing the existing rows 880-977 down to 883-980.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before the current row 880, shifting everything
# below (through row 977) down to rows 883-980.
$ws.Rows("880:882").Insert()

# New date for this week's record: 2023-01-20 (Excel serial 44946)
$newDate = [DateTime]::FromOADate(44946)

# Row 880: Pintón
$ws.Range("A880").Value = 8
$ws.Range("B880").Value = "Terminal La Palmera de La Serena"
$ws.Range("C880").Value = "Coquimbo"
$ws.Range("D880").Value = $newDate
$ws.Range("E880").Value = 4
$ws.Range("F880").Value = "Fruta"
$ws.Range("G880").Value = 100108
$ws.Range("H880").Value = "Tropicales y subtropicales"
$ws.Range("I880").Value = 100108006
$ws.Range("J880").Value = "Plátano"
$ws.Range("K880").Value = "Sin especificar"
$ws.Range("L880").Value = "Pintón"
$ws.Range("M880").Value = 80
$ws.Range("N880").Value = 23000
$ws.Range("O880").Value = 23000
$ws.Range("P880").Value = 23000
$ws.Range("Q880").Value = "$/caja 20 kilos"
$ws.Range("R880").Value = "Ecuador"
$ws.Range("S880").Value = 1150
$ws.Range("T880").Value = 20

# Row 881: Primera Maduro
$ws.Range("A881").Value = 8
$ws.Range("B881").Value = "Terminal La Palmera de La Serena"
$ws.Range("C881").Value = "Coquimbo"
$ws.Range("D881").Value = $newDate
$ws.Range("E881").Value = 4
$ws.Range("F881").Value = "Fruta"
$ws.Range("G881").Value = 100108
$ws.Range("H881").Value = "Tropicales y subtropicales"
$ws.Range("I881").Value = 100108006
$ws.Range("J881").Value = "Plátano"
$ws.Range("K881").Value = "Sin especificar"
$ws.Range("L881").Value = "Primera Maduro"
$ws.Range("M881").Value = 120
$ws.Range("N881").Value = 26000
$ws.Range("O881").Value = 26000
$ws.Range("P881").Value = 26000
$ws.Range("Q881").Value = "$/caja 20 kilos"
$ws.Range("R881").Value = "Ecuador"
$ws.Range("S881").Value = 1300
$ws.Range("T881").Value = 20

# Row 882: Primera Pintón
$ws.Range("A882").Value = 8
$ws.Range("B882").Value = "Terminal La Palmera de La Serena"
$ws.Range("C882").Value = "Coquimbo"
$ws.Range("D882").Value = $newDate
$ws.Range("E882").Value = 4
$ws.Range("F882").Value = "Fruta"
$ws.Range("G882").Value = 100108
$ws.Range("H882").Value = "Tropicales y subtropicales"
$ws.Range("I882").Value = 100108006
$ws.Range("J882").Value = "Plátano"
$ws.Range("K882").Value = "Sin especificar"
$ws.Range("L882").Value = "Primera Pintón"
$ws.Range("M882").Value = 120
$ws.Range("N882").Value = 27000
$ws.Range("O882").Value = 27000
$ws.Range("P882").Value = 27000
$ws.Range("Q882").Value = "$/caja 20 kilos"
$ws.Range("R882").Value = "Ecuador"
$ws.Range("S882").Value = 1350
$ws.Range("T882").Value = 20
